$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 57 reuses the same contact (A/B) as row 55 (seneca62787 / Abdul L Mans),
# adds two new outreach message strings in D/E, and flips the "message sent" flags.
$ws.Range("A57").Value = $ws.Range("A55").Value()
$ws.Range("B57").Value = $ws.Range("B55").Value()
$ws.Range("C57").Value = $true
$ws.Range("D57").Value = "Such a sweet, genuine moment! Capturing these fleeting connections is everything. Do you ever find it challenging to balance documenting & being present? Just sent you something in DMs you might find interesting - check it when you get a chance!"
$ws.Range("E57").Value = "Hey! Just saw your video about the new Supra – awesome content! I'm a big car enthusiast myself and really enjoyed your breakdown of the specs. I actually do some editing and motion graphics for a few YouTube channels in the tech space, helping them level up their visuals.I recently helped a channel jump from 2k to 10k subscribers, mostly by refining their video structure and adding some dynamic transitions. I’m Visuals, by the way. Happy to chat about your content if you're ever looking for a fresh perspective!"
$ws.Range("F57").Value = $true
$ws.Range("G57").Value = $false
